# Add draft "Sex" extension row for version 6.1.0.
#
# The source data had a duplicated "Birth Sex"/"US Core Birth Sex Extension"
# entry across rows 73 and 74. This edit repurposes row 73 into a new
# "Sex" / "US Core Sex Extension" entry, and removes the old duplicate
# row (previously row 74), which shifts every subsequent row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old duplicate "Birth Sex" row; this shifts rows 75-116 up to 74-115.
$ws.Rows.Item(74).Delete()

# Update (now-merged) row 73 to describe the new draft Sex extension.
$ws.Cells.Item(73, 3).Value = "US Core Patient Profile | US Core Sex Extension"
$ws.Cells.Item(73, 4).Value = "Sex"
$ws.Cells.Item(73, 6).Value = "US Core Sex Extension"
